$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.858.38"
$ws.Range("E2").Value = "  +3.20%  "

# Row 3
$ws.Range("D3").Value = "3.521.40"
$ws.Range("E3").Value = "  +2.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.61%  "

# Row 13
$ws.Range("E13").Value = "  -2.69%  "

# Row 14
$ws.Range("D14").Value = "4.077.99"
$ws.Range("E14").Value = "  +2.92%  "

# Row 15
$ws.Range("E15").Value = "  -0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.29%  "

# Row 17
$ws.Range("D17").Value = "3.532.10"
$ws.Range("E17").Value = "  +2.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.14%  "

# Row 20
$ws.Range("D20").Value = "63.705.50"
$ws.Range("E20").Value = "  +2.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.50%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "

# Row 31
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.70%  "

# Row 35
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0494"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.47%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.13%  "

# Row 40
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.20%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.50%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "

# Row 45
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.320"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.75%  "

# Row 47
$ws.Range("D47").Value = "0.0₃0583"
$ws.Range("E47").Value = "  +29.78%  "

# Row 48
$ws.Range("E48").Value = "  +10.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.144"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
